$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: duplicate the last data row-pair (211:212) into new rows 213:214
# (this also extends the sheet dimension from R212 to R214)
$src = $ws.Range("A211:R212")
$dst = $ws.Range("A213:R214")
$src.Copy($dst)

# Step 2: weekly roll-forward of the date-keyed price rows.
# Each row-pair (Primera/Segunda quality) takes the Fecha (D), Volumen (J),
# Precio promedio ponderado (M) and Precio $/Kg (P) of a newly-reported week;
# row 137/138 gets a brand-new week, all others shift down one slot.
# columns: firstRow, Fecha, VolumenPrimera, PrecioPromedioPrimera, PrecioKgPrimera, VolumenSegunda
$updates = @(
    @(137, 44488, 800, 650, 130, 400),
    @(139, 44168, 600, 650, 130, 300),
    @(141, 44292, 600, 650, 130, 300),
    @(143, 44161, 600, 650, 130, 300),
    @(145, 44280, 800, 650, 130, 400),
    @(147, 44274, 600, 650, 130, 300),
    @(149, 44229, 600, 650, 130, 300),
    @(151, 44344, 600, 650, 130, 300),
    @(153, 44358, 600, 650, 130, 300),
    @(155, 44231, 200, 650, 130, 100),
    @(157, 44320, 800, 650, 130, 400),
    @(159, 44474, 600, 650, 130, 300),
    @(161, 44252, 800, 650, 130, 400),
    @(163, 44204, 600, 650, 130, 300),
    @(165, 44362, 200, 650, 130, 100),
    @(167, 44300, 600, 650, 130, 300),
    @(169, 44372, 600, 650, 130, 300),
    @(171, 44350, 600, 650, 130, 300),
    @(173, 44448, 600, 650, 130, 300),
    @(175, 44243, 800, 650, 130, 400),
    @(177, 44385, 600, 650, 130, 300),
    @(179, 44202, 800, 650, 130, 400),
    @(181, 44453, 800, 650, 130, 400),
    @(183, 44435, 1400, 650, 130, 700),
    @(185, 44307, 500, 660, 132, 200),
    @(187, 44159, 600, 650, 130, 300),
    @(189, 44166, 600, 650, 130, 300),
    @(191, 44334, 600, 650, 130, 300),
    @(193, 44476, 100, 650, 130, 50),
    @(195, 44386, 800, 650, 130, 400),
    @(197, 44306, 600, 650, 130, 300),
    @(199, 44357, 800, 650, 130, 400),
    @(201, 44321, 600, 650, 130, 300),
    @(203, 44397, 600, 650, 130, 300),
    @(205, 44314, 400, 650, 130, 200),
    @(207, 44425, 600, 650, 130, 300),
    @(209, 44390, 600, 650, 130, 300),
    @(211, 44250, 600, 650, 130, 300),
    @(213, 44432, 600, 650, 130, 300)
)

foreach ($u in $updates) {
    $rOdd  = $u[0]
    $fecha = $u[1]
    $volP  = $u[2]
    $precioProm = $u[3]
    $precioKg   = $u[4]
    $volS  = $u[5]
    $rEven = $rOdd + 1

    $ws.Range("D$rOdd").Value = $fecha
    $ws.Range("D$rEven").Value = $fecha
    $ws.Range("J$rOdd").Value = $volP
    $ws.Range("J$rEven").Value = $volS
    $ws.Range("M$rOdd").Value = $precioProm
    $ws.Range("P$rOdd").Value = $precioKg
}
